$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Feuil1"

# Fix the "transfert_option" typo -> "transfer_option" in rows 59-70, columns A and B
# Row 65 and 68 also get a semantic change in column B.
$ws.Range("A59").Value = "transfer_option_xlink_role"
$ws.Range("B59").Value = "resources{}.xlink_role"

$ws.Range("A60").Value = "transfer_option_url"
$ws.Range("B60").Value = "resources{}.url"

$ws.Range("A61").Value = "transfer_option_protocol"
$ws.Range("B61").Value = "resources{}.protocol"

$ws.Range("A62").Value = "transfer_option_name"
$ws.Range("B62").Value = "resources{}.name"

$ws.Range("A63").Value = "transfer_option_name_other_lang"
$ws.Range("B63").Value = "resources{}.name_fr"

$ws.Range("A64").Value = "transfer_option_name_other_lang_locale"
$ws.Range("B64").Value = "transfer_option_name_other_lang_locale"

$ws.Range("A65").Value = "transfer_option_description_service"
$ws.Range("B65").Value = "resources{}.transfer_option_description_service"

$ws.Range("A66").Value = "transfer_option_description_format"
$ws.Range("B66").Value = "resources{}.format"

$ws.Range("A67").Value = "transfer_option_description_language"
$ws.Range("B67").Value = "resources{}.transfer_option_description_language"

$ws.Range("A68").Value = "transfer_option_description_service_other_lang"
$ws.Range("B68").Value = "resources{}.transfer_option_description_service_other_lang"

$ws.Range("A69").Value = "transfer_option_description_service_other_lang_locale"
$ws.Range("B69").Value = "service_other_lang_locale"

$ws.Range("A70").Value = "transfer_option_description_language_other_lang"
$ws.Range("B70").Value = "resources{}.transfer_option_description_language_other_lang"

# Update sheet view: scroll position and selection
$ws.Application.ActiveWindow.ScrollRow = 55
$ws.Range("B12").Select()
